# The "2024" sheet tracks bank/service notification text snippets (col R/S)
# and hdfc broadband reminders (col P/Q), each a time-ordered stack with the
# newest entry at the top of its block. A new axis-bank notification and a
# new hdfc notification arrived, so every existing entry in each block
# shifts down one row, the newest text goes on top, and the "Broadband"
# row label (col A) shifts down from row 68 to row 69 to stay aligned with
# its (now shifted) hdfc data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- R/S "axis" notification stack: rows 30-63 shift down by one ---
$ws.Range("R30").Value = "share anyone axis"
$ws.Range("S30").Value = "2024-09-05 16:06:05"
$ws.Range("R31").Value = "internet bal axisbank"
$ws.Range("S31").Value = "2024-09-05 16:05:55"
$ws.Range("R32").Value = "transfer share anyone axis"
$ws.Range("S32").Value = "2024-09-05 16:03:14"
$ws.Range("R33").Value = "axis"
$ws.Range("R34").Value = "your net internet"
$ws.Range("S34").Value = "2024-09-05 15:57:15"
$ws.Range("R35").Value = "hear your feedback atm"
$ws.Range("S35").Value = "2024-09-05 14:21:08"
$ws.Range("S36").Value = "2024-09-05 14:18:32"
$ws.Range("S37").Value = "2024-09-05 14:13:16"
$ws.Range("R38").Value = "axis bna"
$ws.Range("S38").Value = "2024-09-05 14:15:23"
$ws.Range("R39").Value = "balance your axis"
$ws.Range("S39").Value = "2024-09-05 09:20:57"
$ws.Range("R40").Value = "bal axis"
$ws.Range("S40").Value = "2024-09-05 09:06:25"
$ws.Range("R41").Value = "broker"
$ws.Range("S41").Value = "2024-09-04 21:20:47"
$ws.Range("R42").Value = "exclusive on axis"
$ws.Range("S42").Value = "2024-09-04 13:21:05"
$ws.Range("R43").Value = "your corporate axis"
$ws.Range("S43").Value = "2024-09-04 11:46:10"
$ws.Range("R44").Value = "balance your axis"
$ws.Range("S44").Value = "2024-09-04 08:14:16"
$ws.Range("R45").Value = "axis"
$ws.Range("S45").Value = "2024-09-04 07:02:13"
$ws.Range("R46").Value = "bal axisbank w axis"
$ws.Range("S46").Value = "2024-09-04 06:53:15"
$ws.Range("R47").Value = "logging iob internet"
$ws.Range("S47").Value = "2024-09-03 20:09:12"
$ws.Range("R48").Value = "password internet"
$ws.Range("S48").Value = "2024-09-03 20:05:31"
$ws.Range("R49").Value = "logging iob internet"
$ws.Range("S49").Value = "2024-09-03 20:05:09"
$ws.Range("R50").Value = "internet"
$ws.Range("S50").Value = "2024-09-03 19:58:18"
$ws.Range("S51").Value = "2024-09-03 19:54:49"
$ws.Range("R52").Value = "login internet invalid"
$ws.Range("S52").Value = "2024-09-03 19:56:17"
$ws.Range("R53").Value = "corporate internet share"
$ws.Range("S53").Value = "2024-09-03 19:22:58"
$ws.Range("R54").Value = "login sbi internet personal do not share anyone"
$ws.Range("S54").Value = "2024-09-03 19:17:10"
$ws.Range("R55").Value = "login internet personal share"
$ws.Range("S55").Value = "2024-09-03 19:13:40"
$ws.Range("R56").Value = "internet verify it"
$ws.Range("S56").Value = "2024-09-03 19:05:49"
$ws.Range("R57").Value = "balance your axis"
$ws.Range("S57").Value = "2024-09-03 13:14:06"
$ws.Range("R58").Value = "lounge"
$ws.Range("S58").Value = "2024-09-03 13:08:08"
$ws.Range("R59").Value = "balance your axis"
$ws.Range("S59").Value = "2024-09-03 11:21:30"
$ws.Range("R60").Value = "broker"
$ws.Range("S60").Value = "2024-09-01 22:35:38"
$ws.Range("S61").Value = "2024-09-01 10:12:03"
$ws.Range("S62").Value = "2024-09-01 09:42:38"
$ws.Range("S63").Value = "2024-09-01 09:29:24"

# Row 63's old R/S entry ("amazeloan" / 2024-09-01 09:27:06) overflows into
# row 64, which previously held the top of the P/Q "hdfc" stack.
$ws.Range("P64").Value = ""
$ws.Range("Q64").Value = ""
$ws.Range("R64").Value = "amazeloan"
$ws.Range("S64").Value = "2024-09-01 09:27:06"

# --- P/Q "hdfc" notification stack: rows 65-67 shift down by one ---
$ws.Range("Q65").Value = "2024-08-30 12:15:48"
$ws.Range("Q66").Value = "2024-08-21 20:17:10"
$ws.Range("Q67").Value = "2024-08-21 20:16:45"

# Row 67's old entry overflows into row 68, which previously held the
# "Broadband" row label (col A); that label shifts down to row 69.
$ws.Range("A68").Value = ""
$ws.Range("P68").Value = "hdfc"
$ws.Range("Q68").Value = "2024-08-21 20:15:50"
$ws.Range("A69").Value = "Broadband"
